$d = $word.ActiveDocument

# Locate the paragraph that ends the "Requisitos" section
# ("LOQ4031: Química Geral I (Requisito)") - deletion starts right after it,
# so that paragraph itself is left untouched.
$startFind = $d.Content
$startFind.Find.Execute("LOQ4031: Química Geral I (Requisito)", $true, $false,
                         $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$startPos = $startFind.Paragraphs(1).Range.End

# Locate the trailing "powered by Jekyll" copyright paragraph - deletion ends
# at the end of this paragraph (including its paragraph mark), removing it
# along with the "Ver no Jupiter..." paragraph and the blank paragraph
# between them.
$endFind = $d.Content
$endFind.Find.Execute("Creative Commons Attribution", $true, $false,
                       $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$endPos = $endFind.Paragraphs(1).Range.End

$d.Range($startPos, $endPos).Delete()
